$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A136").Value = 135
$ws.Range("B136").Value = "JavaScript - Module"
$ws.Range("C136").Value = "'1801"
$ws.Range("D136").Value = "Create html file for module"
$ws.Range("E136").Interior.Pattern = 17
$ws.Range("E136").Value = "https://www.youtube.com/embed/AjGgauhpK_M"

$ws.Range("A137").Value = 136
$ws.Range("B137").Value = "JavaScript - Module"
$ws.Range("C137").Value = "'1802"
$ws.Range("D137").Value = "Export menu object array"
$ws.Range("E137").Value = "https://www.youtube.com/embed/68ZUeNVGu0k"

$ws.Range("A138").Value = 137
$ws.Range("B138").Value = "JavaScript - Module"
$ws.Range("C138").Value = "'1803"
$ws.Range("D138").Value = "Comment external script in html"
$ws.Range("E138").Value = "https://www.youtube.com/embed/rmn03WMdiaU"

$ws.Range("A139").Value = 138
$ws.Range("B139").Value = "JavaScript - Module"
$ws.Range("C139").Value = "'1804"
$ws.Range("D139").Value = "Convert script to module"
$ws.Range("E139").Value = "https://www.youtube.com/embed/5d_w0JU_FQs"

$ws.Range("A140").Value = 139
$ws.Range("B140").Value = "JavaScript - Module"
$ws.Range("C140").Value = "'1805"
$ws.Range("D140").Value = "Import to use menu object array"
$ws.Range("E140").Value = "https://www.youtube.com/embed/aqO6mycyoLo"

$ws.Range("A141").Value = 140
$ws.Range("B141").Value = "JavaScript - Module"
$ws.Range("C141").Value = "'1806"
$ws.Range("D141").Value = "Advantages of modules explained"
$ws.Range("E141").Value = "https://www.youtube.com/embed/oRxvduGkg7E"

$ws.Range("A142").Value = 141
$ws.Range("B142").Value = "JavaScript - 'if' statement"
$ws.Range("C142").Value = "'1901"
$ws.Range("D142").Value = "Create html for 'if' statement demo"
$ws.Range("E142").Value = "https://www.youtube.com/embed/02hU3q0xOOg"

$ws.Range("A143").Value = 142
$ws.Range("B143").Value = "JavaScript - 'if' statement"
$ws.Range("C143").Value = "'1902"
$ws.Range("D143").Value = "Include 'if' statement"
$ws.Range("E143").Value = "https://www.youtube.com/embed/SaQ3tYC1hWU"

$ws.Range("A144").Value = 143
$ws.Range("B144").Value = "JavaScript - 'if' statement"
$ws.Range("C144").Value = "'1903"
$ws.Range("D144").Value = "'Conditional 'if' statement explained"
$ws.Range("E144").Value = "https://www.youtube.com/embed/KHzocx4_jSg"

$ws.Range("A145").Value = 144
$ws.Range("B145").Value = "JavaScript - 'if' statement"
$ws.Range("C145").Value = "'1904"
$ws.Range("D145").Value = "Code indentation rules"
$ws.Range("E145").Value = "https://www.youtube.com/embed/R1C6ytRg9iY"

$ws.Range("A146").Value = 145
$ws.Range("B146").Value = "JavaScript - Function"
$ws.Range("C146").Value = "'2001"
$ws.Range("D146").Value = "Significance of functions"
$ws.Range("E146").Value = "https://www.youtube.com/embed/u0crhZSNhDo"

$ws.Range("A147").Value = 146
$ws.Range("B147").Value = "JavaScript - Function"
$ws.Range("C147").Value = "'2002"
$ws.Range("D147").Value = "Definition of function"
$ws.Range("E147").Value = "https://www.youtube.com/embed/t8w6FQWcNac"

$ws.Range("A148").Value = 147
$ws.Range("B148").Value = "JavaScript - Function"
$ws.Range("C148").Value = "'2003"
$ws.Range("D148").Value = "Calling the function"
$ws.Range("E148").Value = "https://www.youtube.com/embed/5KKgepsCUyA"

$ws.Range("A149").Value = 148
$ws.Range("B149").Value = "JavaScript - Function"
$ws.Range("C149").Value = "'2004"
$ws.Range("D149").Value = "Function and code execution flow"
$ws.Range("E149").Value = "https://www.youtube.com/embed/znnr8G_vwKU"

$ws.Range("A150").Value = 149
$ws.Range("B150").Value = "JavaScript - Function"
$ws.Range("C150").Value = "'2005"
$ws.Range("D150").Value = "Passing values to a function"
$ws.Range("E150").Value = "https://www.youtube.com/embed/AVVV1uG3Fn0"

$ws.Range("A151").Value = 150
$ws.Range("B151").Value = "JavaScript - Function"
$ws.Range("C151").Value = "'2006"
$ws.Range("D151").Value = "Return value from a function"
$ws.Range("E151").Value = "https://www.youtube.com/embed/VB_rokYGjlk"

$ws.Range("A152").Value = 151
$ws.Range("B152").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C152").Value = "'2101"
$ws.Range("D152").Value = "Displaying all menu items"
$ws.Range("E152").Value = "https://www.youtube.com/embed/Bw2Zk-Yh7yI"

$ws.Range("A153").Value = 152
$ws.Range("B153").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C153").Value = "'2102"
$ws.Range("D153").Value = "Select all image files"
$ws.Range("E153").Value = "https://www.youtube.com/embed/GRxP4AxlcNA"

$ws.Range("A154").Value = 153
$ws.Range("B154").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C154").Value = "'2103"
$ws.Range("D154").Value = "Exclude first menu item in selection and copy"
$ws.Range("E154").Value = "https://www.youtube.com/embed/5NQLwYxAlhk"

$ws.Range("A155").Value = 154
$ws.Range("B155").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C155").Value = "'2104"
$ws.Range("D155").Value = "Copy menu object array script file"
$ws.Range("E155").Value = "https://www.youtube.com/embed/XaItWNkMLb0"

$ws.Range("A156").Value = 155
$ws.Range("B156").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C156").Value = "'2105"
$ws.Range("D156").Value = "Create script for menu page and include import"
$ws.Range("E156").Value = "https://www.youtube.com/embed/D1rAdOAlCJU"

$ws.Range("A157").Value = 156
$ws.Range("B157").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C157").Value = "'2106"
$ws.Range("D157").Value = "Logic to display menu items"
$ws.Range("E157").Value = "https://www.youtube.com/embed/u4OUHLA3u-o"

$ws.Range("A158").Value = 157
$ws.Range("B158").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C158").Value = "'2107"
$ws.Range("D158").Value = "Menu page html code formatting"
$ws.Range("E158").Value = "https://www.youtube.com/embed/31N22kvRjA4"

$ws.Range("A159").Value = 158
$ws.Range("B159").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C159").Value = "'2108"
$ws.Range("D159").Value = "Collapse menu item div tag and copy"
$ws.Range("E159").Value = "https://www.youtube.com/embed/Ga2eNCsjQMw"

$ws.Range("A160").Value = 159
$ws.Range("B160").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C160").Value = "'2109"
$ws.Range("D160").Value = "Define div tag template literal in menu script"
$ws.Range("E160").Value = "https://www.youtube.com/embed/E58jNx7m5qc"

$ws.Range("A161").Value = 160
$ws.Range("B161").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C161").Value = "'2110"
$ws.Range("D161").Value = "Include template literal placeholders"
$ws.Range("E161").Value = "https://www.youtube.com/embed/rby7LqIAw-4"

$ws.Range("A162").Value = 161
$ws.Range("B162").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C162").Value = "'2111"
$ws.Range("D162").Value = "Include properties in placeholders"
$ws.Range("E162").Value = "https://www.youtube.com/embed/OxdZrX6Yg4U"

$ws.Range("A163").Value = 162
$ws.Range("B163").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C163").Value = "'2112"
$ws.Range("D163").Value = "Include 'id' for row div"
$ws.Range("E163").Value = "https://www.youtube.com/embed/1hf3mYLWdXY"

$ws.Range("A164").Value = 163
$ws.Range("B164").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C164").Value = "'2113"
$ws.Range("D164").Value = "Comment row div content"
$ws.Range("E164").Value = "https://www.youtube.com/embed/3y_DzfrdgEc"

$ws.Range("A165").Value = 164
$ws.Range("B165").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C165").Value = "'2114"
$ws.Range("D165").Value = "Include menu script in menu page"
$ws.Range("E165").Value = "https://www.youtube.com/embed/h2kZVL9nQh8"

$ws.Range("A166").Value = 165
$ws.Range("B166").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C166").Value = "'2115"
$ws.Range("D166").Value = "Display menu details using DOM"
$ws.Range("E166").Value = "https://www.youtube.com/embed/hwKMQsYvh6c"

$ws.Range("A167").Value = 166
$ws.Range("B167").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C167").Value = "'2116"
$ws.Range("D167").Value = "Display only active items using 'if' condition"
$ws.Range("E167").Value = "https://www.youtube.com/embed/giQG9dQtJtw"

$ws.Range("A168").Value = 167
$ws.Range("B168").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C168").Value = "'2117"
$ws.Range("D168").Value = "Function to convert menu object to html"
$ws.Range("E168").Value = "https://www.youtube.com/embed/eIu46oT_TgI"

$ws.Range("A169").Value = 168
$ws.Range("B169").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C169").Value = "'2118"
$ws.Range("D169").Value = "Advantage of function of menu div conversion"
$ws.Range("E169").Value = "https://www.youtube.com/embed/EYhdc7T6BSY"

$ws.Range("A170").Value = 169
$ws.Range("B170").Value = "Patisserie - Display menu items using JavaScript menu objects"
$ws.Range("C170").Value = "'2119"
$ws.Range("D170").Value = "Bottom margin for menu items"
$ws.Range("E170").Value = "https://www.youtube.com/embed/q5ce1nwokCk"

$ws.Range("E137").Select()
